$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.216.76"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.864.90"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7177"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.55"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07752"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3075"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.03"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08253"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "1.886.40"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7157"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.211"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.30"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "29.212.67"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.04"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007789"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Value = "2.121.45"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.12"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.925"
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1593"
$ws.Range("E25").Value = "  +8.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.12"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.908"
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.17"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.492"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.310"
$ws.Range("E30").Value = "  -4.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.359"
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.082"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05195"
$ws.Range("E33").Value = "  -0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.916"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7278"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.691"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "1.162.88"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9022"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.093"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.16"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "2.016.99"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.763"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.260"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9980"
$ws.Range("E51").Value = "  -0.52%  "
